# Weekly update: a new week's price record is inserted at the top of the
# data (row 170), pushing all subsequent rows down by one. The new record
# is the same "Ajo - Chino - Primera" quote that was most recently seen
# (originally at row 199, dated 2021-11-08 / serial 44508) but reported two
# days later (2021-11-10 / serial 44510). The last row of the sheet
# (originally row 248) ends up at row 249, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 170; rows 170..248 shift down to 171..249.
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new weekly record.
$ws.Range("A170").Value = 3
$ws.Range("B170").Value = "Femacal de La Calera"
$ws.Range("C170").Value = "Coquimbo"
$ws.Range("D170").Value = 44510
$ws.Range("E170").Value = 5
$ws.Range("F170").Value = 100112003
$ws.Range("G170").Value = "Ajo"
$ws.Range("H170").Value = "Chino"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 73
$ws.Range("K170").Value = 16000
$ws.Range("L170").Value = 16500
$ws.Range("M170").Value = 16260
$ws.Range("N170").Value = "`$/caja 10 kilos"
$ws.Range("O170").Value = "China"
$ws.Range("P170").Value = 1626
$ws.Range("Q170").Value = 10
$ws.Range("R170").Value = "Hortaliza"
